# Insert a new data row before existing row 855, shifting subsequent rows down.
# Row 854 already has the same date ("2026/02/25") stored as a plain text
# value, so copy it (to inherit the text-typed date/weekday cells without
# Excel's automatic date-literal conversion or any new cell style), insert
# the copy above row 855 (pushing old rows 855+ down by one), then overwrite
# the time/ranking columns with the new row's actual values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(854).Copy()
$ws.Rows.Item(855).Insert()

$ws.Range("C855").Value = 13
$ws.Range("D855").Value = 29
